# Update the "想去人数" (want-to-go count) figures in column F across the
# four worksheets of the workbook to the newly scraped values.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        [string]$SheetName,
        [hashtable]$RowValues
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Range("F$row").Value = $RowValues[$row]
    }
}

# Sheet "展览" (Exhibitions)
Set-FValues "展览" @{
    2  = 71
    6  = 333
    9  = 70
    10 = 33
    11 = 673
    12 = 1513
    13 = 5868
    15 = 1654
    17 = 5583
    18 = 101
    22 = 1585
    25 = 97
    26 = 1180
    27 = 683
    31 = 3837
}

# Sheet "演出" (Performances)
Set-FValues "演出" @{
    4 = 106
    5 = 207
    8 = 318
    9 = 2
}

# Sheet "本地生活" (Local Life)
Set-FValues "本地生活" @{
    2 = 9458
    5 = 550
}

# Sheet "全部类型" (All Types)
Set-FValues "全部类型" @{
    2  = 9458
    5  = 71
    7  = 550
    8  = 333
    11 = 70
    14 = 673
    15 = 1513
    16 = 5868
    18 = 318
    19 = 1654
    20 = 2
    25 = 5584
    26 = 101
    30 = 1585
    33 = 97
    34 = 1180
    35 = 683
    45 = 3837
}
